$wb = $excel.ActiveWorkbook

# --- krasnoludAncestry sheet: update B-column numeric values ---
$ws1 = $wb.Worksheets.Item("krasnoludAncestry")

$ws1.Range("B2").Value = 10
$ws1.Range("B3").Value = 9
$ws1.Range("B5").Value = 10
$ws1.Range("B6").Value = 11
$ws1.Range("B7").Value = 14
$ws1.Range("B8").Value = 8
$ws1.Range("B9").Value = 9
$ws1.Range("B10").Value = 3

# C15 currently holds the language choice "Elifcki"; replace it with the
# new "Krasnoludzki" option, and add the same choice again in C16 (new
# quirks/language choice row under "Jezyk Pismo:").
$ws1.Range("C15").Value = "Krasnoludzki"
$ws1.Range("C16").Value = "Krasnoludzki"

# Move the active selection on this sheet and make it the active tab.
$ws1.Activate() | Out-Null
$ws1.Range("E29").Select() | Out-Null
